# Applies the weekly Fruit/Vegetable price update for the
# "Macroferia Regional de Talca - Granada" sheet: existing rows 2-15
# get their Fecha/Calidad/Volumen/Precio.../Origen/Kg-unidad fields
# reshuffled onto different dates (per the source diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44354
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 18000
$ws.Range("S2").Value = 1000

# Row 3
$ws.Range("D3").Value = 44355
$ws.Range("L3").Value = 'Especial'
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("R3").Value = 'Provincia de Limarí'
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 44342
$ws.Range("L4").Value = 'Especial'
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("S4").Value = 1111

# Row 5
$ws.Range("D5").Value = 44340
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 230

# Row 6
$ws.Range("D6").Value = 44319
$ws.Range("L6").Value = 'Especial'
$ws.Range("M6").Value = 120
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 20000
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = 'Provincia de Limarí'
$ws.Range("S6").Value = 1111
$ws.Range("T6").Value = 18

# Row 7
$ws.Range("D7").Value = 44326

# Row 8
$ws.Range("D8").Value = 44316
$ws.Range("M8").Value = 300

# Row 9
$ws.Range("D9").Value = 44294
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("Q9").Value = '$/caja 15 kilos granel'
$ws.Range("R9").Value = 'Región Metropolitana'
$ws.Range("S9").Value = 800
$ws.Range("T9").Value = 15

# Row 10
$ws.Range("D10").Value = 44348
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 20000
$ws.Range("P10").Value = 20000
$ws.Range("S10").Value = 1111

# Row 11
$ws.Range("D11").Value = 44299
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 15000
$ws.Range("Q11").Value = '$/caja 15 kilos granel'
$ws.Range("R11").Value = 'Provincia de Curicó'
$ws.Range("S11").Value = 1000
$ws.Range("T11").Value = 15

# Row 12
$ws.Range("D12").Value = 44358
$ws.Range("L12").Value = 'Especial'
$ws.Range("M12").Value = 150
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 18000
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("R12").Value = 'Provincia de Limarí'
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 44358
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 17000
$ws.Range("O13").Value = 17000
$ws.Range("P13").Value = 17000
$ws.Range("S13").Value = 944

# Row 14
$ws.Range("D14").Value = 44291
$ws.Range("L14").Value = 'Primera'
$ws.Range("N14").Value = 12000
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 12000
$ws.Range("Q14").Value = '$/caja 15 kilos granel'
$ws.Range("R14").Value = 'Región Metropolitana'
$ws.Range("S14").Value = 800
$ws.Range("T14").Value = 15

# Row 15
$ws.Range("D15").Value = 44328
$ws.Range("L15").Value = 'Especial'
$ws.Range("M15").Value = 250
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 20000
$ws.Range("S15").Value = 1111

